# Scope 1 stationary fuel workbook update
# - CH4 / N2O emission factor units changed from kg to g (values scaled x1000)
# - Selection moved to B12
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update CH4 Factor header (D3): "CH4 Factor (kg/ mmBtu)" -> "CH4 Factor (g/ mmBtu)"
# Preserve rich text formatting: "CH" (default/cell font), "4" (bold Arial 10 subscript),
# " Factor (g/ mmBtu)" (bold Arial 10)
$d3 = $ws.Range("D3")
$d3.Value = "CH4 Factor (g/ mmBtu)"

$d3Sub = $d3.Characters(3, 1)
$d3Sub.Font.Bold = $true
$d3Sub.Font.Size = 10
$d3Sub.Font.Name = "Arial"
$d3Sub.Font.Subscript = $true

$d3Suffix = $d3.Characters(4, 18)
$d3Suffix.Font.Bold = $true
$d3Suffix.Font.Size = 10
$d3Suffix.Font.Name = "Arial"

# --- Update N2O Factor header (E3): "N2O Factor (kg / mmBtu)" -> "N2O Factor (g / mmBtu)"
# Preserve rich text formatting: "N" (default/cell font), "2" (bold Arial 10 subscript),
# "O Factor (g / mmBtu)" (bold Arial 10)
$e3 = $ws.Range("E3")
$e3.Value = "N2O Factor (g / mmBtu)"

$e3Sub = $e3.Characters(2, 1)
$e3Sub.Font.Bold = $true
$e3Sub.Font.Size = 10
$e3Sub.Font.Name = "Arial"
$e3Sub.Font.Subscript = $true

$e3Suffix = $e3.Characters(3, 20)
$e3Suffix.Font.Bold = $true
$e3Suffix.Font.Size = 10
$e3Suffix.Font.Name = "Arial"

# --- Rescale CH4 / N2O factor values for Propane and Natural Gas rows (kg -> g, x1000)
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 0.6
$ws.Range("D6").Value = 3
$ws.Range("E6").Value = 0.6

# --- Update selection to B12
$ws.Range("B12").Select()
